$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1083.7894
$ws.Cells.Item(40, 9).Value = 1078.091
$ws.Cells.Item(40, 10).Value = 1091.625
$ws.Cells.Item(40, 11).Value = 1078.091
$ws.Cells.Item(40, 12).Value = 1091.625
$ws.Cells.Item(40, 13).Value = -903.0909999999999
$ws.Cells.Item(40, 14).Value = -1441.625

$ws.Cells.Item(98, 8).Value = 693.1818
$ws.Cells.Item(98, 9).Value = 672.3
$ws.Cells.Item(98, 10).Value = 902
$ws.Cells.Item(98, 11).Value = 672.3
$ws.Cells.Item(98, 12).Value = 902
$ws.Cells.Item(98, 13).Value = 825.7
$ws.Cells.Item(98, 14).Value = -3898

$ws.Cells.Item(112, 8).Value = 1009.67645
$ws.Cells.Item(112, 10).Value = 1041.5312
$ws.Cells.Item(112, 12).Value = 3124.5936
$ws.Cells.Item(112, 14).Value = -5340.5936

$ws.Cells.Item(122, 8).Value = 693.1818
$ws.Cells.Item(122, 9).Value = 672.3
$ws.Cells.Item(122, 10).Value = 902
$ws.Cells.Item(122, 11).Value = 2016.9
$ws.Cells.Item(122, 12).Value = 2706
$ws.Cells.Item(122, 13).Value = 433.1000000000001
$ws.Cells.Item(122, 14).Value = -7606

$ws.Cells.Item(132, 8).Value = 2953.923
$ws.Cells.Item(132, 9).Value = 3074.8333
$ws.Cells.Item(132, 10).Value = 1503
$ws.Cells.Item(132, 11).Value = 9224.499899999999
$ws.Cells.Item(132, 12).Value = 4509
$ws.Cells.Item(132, 13).Value = -6694.499899999999
$ws.Cells.Item(132, 14).Value = -9569

$ws.Cells.Item(137, 8).Value = 1990.9131
$ws.Cells.Item(137, 9).Value = 1883.3889
$ws.Cells.Item(137, 10).Value = 2378
$ws.Cells.Item(137, 11).Value = 5650.1667
$ws.Cells.Item(137, 12).Value = 7134
$ws.Cells.Item(137, 13).Value = -3100.1667
$ws.Cells.Item(137, 14).Value = -12234

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1885.3818
$ws.Cells.Item(32, 9).Value = 1628.8125
$ws.Cells.Item(32, 10).Value = 3644.7144
$ws.Cells.Item(32, 11).Value = 1628.8125
$ws.Cells.Item(32, 12).Value = 3644.7144
$ws.Cells.Item(32, 13).Value = -1341.8125
$ws.Cells.Item(32, 14).Value = -4218.7144

$ws.Cells.Item(45, 8).Value = 3497.1538
$ws.Cells.Item(45, 9).Value = 3821.2
$ws.Cells.Item(45, 10).Value = 3294.625
$ws.Cells.Item(45, 11).Value = 3821.2
$ws.Cells.Item(45, 12).Value = 3294.625
$ws.Cells.Item(45, 13).Value = -3444.2
$ws.Cells.Item(45, 14).Value = -4048.625

$ws.Cells.Item(74, 8).Value = 2320
$ws.Cells.Item(74, 9).Value = 2011.579
$ws.Cells.Item(74, 10).Value = 5250
$ws.Cells.Item(74, 11).Value = 2011.579
$ws.Cells.Item(74, 12).Value = 5250
$ws.Cells.Item(74, 13).Value = -1137.579
$ws.Cells.Item(74, 14).Value = -6998

$ws.Cells.Item(77, 8).Value = 2320
$ws.Cells.Item(77, 9).Value = 2011.579
$ws.Cells.Item(77, 10).Value = 5250
$ws.Cells.Item(77, 11).Value = 10057.895
$ws.Cells.Item(77, 12).Value = 26250
$ws.Cells.Item(77, 13).Value = -5689.895
$ws.Cells.Item(77, 14).Value = -34986

$ws.Cells.Item(88, 8).Value = 87037.836
$ws.Cells.Item(88, 10).Value = 104095.4
$ws.Cells.Item(88, 12).Value = 104095.4
$ws.Cells.Item(88, 14).Value = -104907.4

$ws.Cells.Item(91, 8).Value = 87037.836
$ws.Cells.Item(91, 10).Value = 104095.4
$ws.Cells.Item(91, 12).Value = 104095.4
$ws.Cells.Item(91, 14).Value = -106903.4

$ws.Cells.Item(97, 8).Value = 1603.4736
$ws.Cells.Item(97, 9).Value = 1343.1875
$ws.Cells.Item(97, 10).Value = 2991.6667
$ws.Cells.Item(97, 11).Value = 1343.1875
$ws.Cells.Item(97, 12).Value = 2991.6667
$ws.Cells.Item(97, 13).Value = -847.1875
$ws.Cells.Item(97, 14).Value = -3983.6667

$ws.Cells.Item(110, 8).Value = 3155.4546
$ws.Cells.Item(110, 9).Value = 3230
$ws.Cells.Item(110, 10).Value = 3025
$ws.Cells.Item(110, 11).Value = 3230
$ws.Cells.Item(110, 12).Value = 3025
$ws.Cells.Item(110, 13).Value = -1185
$ws.Cells.Item(110, 14).Value = -7115

$ws.Cells.Item(132, 8).Value = 34461.438
$ws.Cells.Item(132, 9).Value = 2919.2222
$ws.Cells.Item(132, 10).Value = 75015.71000000001
$ws.Cells.Item(132, 11).Value = 8757.6666
$ws.Cells.Item(132, 12).Value = 225047.13
$ws.Cells.Item(132, 13).Value = -6227.6666
$ws.Cells.Item(132, 14).Value = -230107.13

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2626.4194
$ws.Cells.Item(94, 9).Value = 1168.1666
$ws.Cells.Item(94, 10).Value = 4645.5386
$ws.Cells.Item(94, 11).Value = 1168.1666
$ws.Cells.Item(94, 12).Value = 4645.5386
$ws.Cells.Item(94, 13).Value = -717.1666
$ws.Cells.Item(94, 14).Value = -5547.5386

$ws.Cells.Item(105, 8).Value = 3850
$ws.Cells.Item(105, 9).Value = 3266.6667
$ws.Cells.Item(105, 10).Value = 5600
$ws.Cells.Item(105, 11).Value = 3266.6667
$ws.Cells.Item(105, 12).Value = 5600
$ws.Cells.Item(105, 13).Value = -1519.6667
$ws.Cells.Item(105, 14).Value = -9094

$ws.Cells.Item(107, 8).Value = 1898.8334
$ws.Cells.Item(107, 9).Value = 845
$ws.Cells.Item(107, 11).Value = 845
$ws.Cells.Item(107, 13).Value = 1075

$ws.Cells.Item(128, 8).Value = 3000
$ws.Cells.Item(128, 9).Value = 3000
$ws.Cells.Item(128, 11).Value = 9000
$ws.Cells.Item(128, 13).Value = -6510

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10579.381
$ws.Cells.Item(31, 9).Value = 12481.363
$ws.Cells.Item(31, 10).Value = 3605.4443
$ws.Cells.Item(31, 11).Value = 12481.363
$ws.Cells.Item(31, 12).Value = 3605.4443
$ws.Cells.Item(31, 13).Value = -12186.363
$ws.Cells.Item(31, 14).Value = -4195.4443

$ws.Cells.Item(34, 8).Value = 10579.381
$ws.Cells.Item(34, 9).Value = 12481.363
$ws.Cells.Item(34, 10).Value = 3605.4443
$ws.Cells.Item(34, 11).Value = 12481.363
$ws.Cells.Item(34, 12).Value = 3605.4443
$ws.Cells.Item(34, 13).Value = -12279.363
$ws.Cells.Item(34, 14).Value = -4009.4443

$ws.Cells.Item(132, 8).Value = 26551.096
$ws.Cells.Item(132, 9).Value = 36910.645
$ws.Cells.Item(132, 10).Value = 5832
$ws.Cells.Item(132, 11).Value = 110731.935
$ws.Cells.Item(132, 12).Value = 17496
$ws.Cells.Item(132, 13).Value = -108201.935
$ws.Cells.Item(132, 14).Value = -22556

$ws.Cells.Item(135, 8).Value = 50480
$ws.Cells.Item(135, 10).Value = 50480
$ws.Cells.Item(135, 12).Value = 50480
$ws.Cells.Item(135, 14).Value = -60620

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 21.5
$ws.Cells.Item(2, 10).Value = 30.666666
$ws.Cells.Item(2, 12).Value = 183.999996
$ws.Cells.Item(2, 14).Value = -409.999996

$ws.Cells.Item(5, 8).Value = 1265.8334
$ws.Cells.Item(5, 9).Value = 1018.5
$ws.Cells.Item(5, 11).Value = 3055.5
$ws.Cells.Item(5, 13).Value = -2943.5

$ws.Cells.Item(22, 8).Value = 20870
$ws.Cells.Item(22, 9).Value = 50550
$ws.Cells.Item(22, 10).Value = 1083.3334
$ws.Cells.Item(22, 11).Value = 151650
$ws.Cells.Item(22, 12).Value = 3250.0002
$ws.Cells.Item(22, 13).Value = -151481
$ws.Cells.Item(22, 14).Value = -3588.0002

$ws.Cells.Item(27, 8).Value = 20870
$ws.Cells.Item(27, 9).Value = 50550
$ws.Cells.Item(27, 10).Value = 1083.3334
$ws.Cells.Item(27, 11).Value = 151650
$ws.Cells.Item(27, 12).Value = 3250.0002
$ws.Cells.Item(27, 13).Value = -151548
$ws.Cells.Item(27, 14).Value = -3454.0002

$ws.Cells.Item(32, 8).Value = 1650
$ws.Cells.Item(32, 9).Value = 1300
$ws.Cells.Item(32, 10).Value = 2000
$ws.Cells.Item(32, 11).Value = 3900
$ws.Cells.Item(32, 12).Value = 6000
$ws.Cells.Item(32, 13).Value = -3617
$ws.Cells.Item(32, 14).Value = -6566

$ws.Cells.Item(39, 8).Value = 2871.5
$ws.Cells.Item(39, 10).Value = 3345.3333
$ws.Cells.Item(39, 12).Value = 10035.9999
$ws.Cells.Item(39, 14).Value = -10623.9999

$ws.Cells.Item(46, 8).Value = 1253.2
$ws.Cells.Item(46, 9).Value = 500
$ws.Cells.Item(46, 10).Value = 1336.8889
$ws.Cells.Item(46, 11).Value = 1500
$ws.Cells.Item(46, 12).Value = 4010.6667
$ws.Cells.Item(46, 13).Value = -1409
$ws.Cells.Item(46, 14).Value = -4192.6667

$ws.Cells.Item(68, 8).Value = 25700.75
$ws.Cells.Item(68, 10).Value = 34001
$ws.Cells.Item(68, 12).Value = 102003
$ws.Cells.Item(68, 14).Value = -103625

$ws.Cells.Item(71, 8).Value = 25700.75
$ws.Cells.Item(71, 10).Value = 34001
$ws.Cells.Item(71, 12).Value = 306009
$ws.Cells.Item(71, 14).Value = -314121

$ws.Cells.Item(122, 8).Value = 636.8
$ws.Cells.Item(122, 10).Value = 636.8
$ws.Cells.Item(122, 12).Value = 5731.2
$ws.Cells.Item(122, 14).Value = -10631.2

$ws.Cells.Item(131, 8).Value = 789.86
$ws.Cells.Item(131, 9).Value = 200
$ws.Cells.Item(131, 10).Value = 795.8182
$ws.Cells.Item(131, 11).Value = 600
$ws.Cells.Item(131, 12).Value = 2387.4546
$ws.Cells.Item(131, 13).Value = 4440
$ws.Cells.Item(131, 14).Value = -12467.4546

$ws.Cells.Item(135, 8).Value = 1265.8334
$ws.Cells.Item(135, 9).Value = 1018.5
$ws.Cells.Item(135, 11).Value = 9166.5
$ws.Cells.Item(135, 13).Value = -6631.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1881.9048
$ws.Cells.Item(97, 9).Value = 1066.7059
$ws.Cells.Item(97, 10).Value = 5346.5
$ws.Cells.Item(97, 11).Value = 1066.7059
$ws.Cells.Item(97, 12).Value = 5346.5
$ws.Cells.Item(97, 13).Value = -570.7058999999999
$ws.Cells.Item(97, 14).Value = -6338.5

$ws.Cells.Item(113, 8).Value = 2679.5
$ws.Cells.Item(113, 9).Value = 1961.3
$ws.Cells.Item(113, 11).Value = 1961.3
$ws.Cells.Item(113, 13).Value = 208.7

$ws.Cells.Item(122, 8).Value = 1922.3
$ws.Cells.Item(122, 9).Value = 1969.2222
$ws.Cells.Item(122, 10).Value = 1500
$ws.Cells.Item(122, 11).Value = 5907.6666
$ws.Cells.Item(122, 12).Value = 4500
$ws.Cells.Item(122, 13).Value = -3457.6666
$ws.Cells.Item(122, 14).Value = -9400

$ws.Cells.Item(132, 8).Value = 21354
$ws.Cells.Item(132, 9).Value = 4157.1177
$ws.Cells.Item(132, 10).Value = 47931
$ws.Cells.Item(132, 11).Value = 12471.3531
$ws.Cells.Item(132, 12).Value = 143793
$ws.Cells.Item(132, 13).Value = -9941.3531
$ws.Cells.Item(132, 14).Value = -148853

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 19500
$ws.Cells.Item(69, 10).Value = 19500
$ws.Cells.Item(69, 12).Value = 19500
$ws.Cells.Item(69, 14).Value = -20998

$ws.Cells.Item(72, 8).Value = 19500
$ws.Cells.Item(72, 10).Value = 19500
$ws.Cells.Item(72, 12).Value = 58500
$ws.Cells.Item(72, 14).Value = -65988

$ws.Cells.Item(81, 8).Value = 2233.3333
$ws.Cells.Item(81, 9).Value = 2180
$ws.Cells.Item(81, 10).Value = 2500
$ws.Cells.Item(81, 11).Value = 4360
$ws.Cells.Item(81, 12).Value = 5000
$ws.Cells.Item(81, 13).Value = -3299
$ws.Cells.Item(81, 14).Value = -7122

$ws.Cells.Item(84, 8).Value = 2233.3333
$ws.Cells.Item(84, 9).Value = 2180
$ws.Cells.Item(84, 10).Value = 2500
$ws.Cells.Item(84, 11).Value = 21800
$ws.Cells.Item(84, 12).Value = 25000
$ws.Cells.Item(84, 13).Value = -16496
$ws.Cells.Item(84, 14).Value = -35608

$ws.Cells.Item(96, 8).Value = 4942.857
$ws.Cells.Item(96, 9).Value = 700
$ws.Cells.Item(96, 10).Value = 5650
$ws.Cells.Item(96, 11).Value = 700
$ws.Cells.Item(96, 12).Value = 5650
$ws.Cells.Item(96, 13).Value = 673
$ws.Cells.Item(96, 14).Value = -8396

$ws.Cells.Item(132, 8).Value = 3522.7693
$ws.Cells.Item(132, 9).Value = 3088.889
$ws.Cells.Item(132, 11).Value = 9266.667000000001
$ws.Cells.Item(132, 13).Value = -6736.667000000001

$ws.Cells.Item(138, 8).Value = 100379
$ws.Cells.Item(138, 10).Value = 100379
$ws.Cells.Item(138, 12).Value = 100379
$ws.Cells.Item(138, 14).Value = -110659
